$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7303
$ws1.Range("F4").Value = 3514
$ws1.Range("F6").Value = 3848
$ws1.Range("F8").Value = 79
$ws1.Range("F9").Value = 79
$ws1.Range("F11").Value = 149
$ws1.Range("F12").Value = 507
$ws1.Range("F15").Value = 368
$ws1.Range("F17").Value = 7
$ws1.Range("F19").Value = 4120
$ws1.Range("F21").Value = 410
$ws1.Range("F24").Value = 1659
$ws1.Range("F25").Value = 113
$ws1.Range("F27").Value = 3033
$ws1.Range("F28").Value = 2225
$ws1.Range("F29").Value = 61
$ws1.Range("F30").Value = 79
$ws1.Range("F32").Value = 34
$ws1.Range("F33").Value = 99
$ws1.Range("F36").Value = 4306
$ws1.Range("F37").Value = 477
$ws1.Range("F38").Value = 322
$ws1.Range("F39").Value = 57
$ws1.Range("F41").Value = 804
$ws1.Range("F42").Value = 210
$ws1.Range("F43").Value = 11
$ws1.Range("F44").Value = 1634
$ws1.Range("F46").Value = 32
$ws1.Range("F47").Value = 606
$ws1.Range("F48").Value = 720

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 439
$ws2.Range("F6").Value = 66
$ws2.Range("F10").Value = 45
$ws2.Range("F16").Value = 578

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 165

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 165
$ws4.Range("F4").Value = 7303
$ws4.Range("F5").Value = 3514
$ws4.Range("F6").Value = 3514
$ws4.Range("F7").Value = 3848
$ws4.Range("F9").Value = 79
$ws4.Range("F10").Value = 79
$ws4.Range("F13").Value = 149
$ws4.Range("F14").Value = 507
$ws4.Range("F15").Value = 66
$ws4.Range("F17").Value = 368
$ws4.Range("F19").Value = 7
$ws4.Range("F21").Value = 4120
$ws4.Range("F23").Value = 45
$ws4.Range("F25").Value = 410
$ws4.Range("F28").Value = 1659
$ws4.Range("F29").Value = 113
$ws4.Range("F31").Value = 3033
$ws4.Range("F32").Value = 2225
$ws4.Range("F33").Value = 61
$ws4.Range("F34").Value = 79
$ws4.Range("F39").Value = 4306
$ws4.Range("F41").Value = 477
$ws4.Range("F42").Value = 322
$ws4.Range("F44").Value = 804
$ws4.Range("F45").Value = 210
$ws4.Range("F46").Value = 1634
$ws4.Range("F48").Value = 32
$ws4.Range("F49").Value = 606
$ws4.Range("F50").Value = 720
